$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.824.01'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '2.452.65'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("D9").Value = '2.456.02'
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0980'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.28%  '
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.338'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.80%  '
$ws.Range("D14").Value = '2.886.84'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '57.744.35'
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000134'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("D18").Value = '2.451.07'
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '319.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.70%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.407'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.06%  '
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("D29").Value = '0.0₃0736'
$ws.Range("E29").Value = '  -4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.62'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.25%  '
$ws.Range("E31").Value = '  -4.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("E37").Value = '  -6.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("E39").Value = '  -4.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '269.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.588'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.44%  '
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0484'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.91%  '
$ws.Range("E49").Value = '  -4.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.47%  '
$ws.Range("D51").Value = '1.721.13'
$ws.Range("E51").Value = '  -2.03%  '
